$wb2 = $excel.ActiveWorkbook
$ws = $wb2.ActiveSheet

# Insert a new row at position 219; this shifts existing rows 219-310 down to 220-311
$ws.Rows(219).Insert()

# Populate the newly inserted row 219 with the new data record
$ws.Range("A219").Value = 4
$ws.Range("B219").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C219").Value = "Los Lagos"
$ws.Range("D219").Value = 44875
$ws.Range("D219").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E219").Value = 10
$ws.Range("F219").Value = 100112032
$ws.Range("G219").Value = "Zapallo italiano"
$ws.Range("H219").Value = "Sin especificar"
$ws.Range("I219").Value = "Primera"
$ws.Range("J219").Value = 180
$ws.Range("K219").Value = 11000
$ws.Range("L219").Value = 11000
$ws.Range("M219").Value = 11000
$ws.Range("N219").Value = "$/caja 50 unidades"
$ws.Range("O219").Value = "Región de O'Higgins"
$ws.Range("P219").Value = 220
$ws.Range("Q219").Value = 50
$ws.Range("R219").Value = "Hortaliza"
